$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Rename "EmailSubscriptionTypes" -> "EmailTypes" and rework its table to
#    add the new "IsSubscription" column + a 3rd row (Email Verification).
# ---------------------------------------------------------------------------
$emailTypes = $wb.Worksheets.Item("EmailSubscriptionTypes")
$emailTypes.Name = "EmailTypes"

# Headers (row 1)
$emailTypes.Range("A1").Value = "EmaiTypeID"
$emailTypes.Range("B1").Value = "EmailType"
$emailTypes.Range("C1").Value = "Description"
$emailTypes.Range("D1").Value = "EmailTypeTidy"
$emailTypes.Range("E1").Value = "IsSubscription"
$emailTypes.Range("F1").Value = "INSERT INTO EmailTypes (EmailTypeID,EmailType,Description,EmailTypeTidy,IsSubscription) VALUES"
$emailTypes.Range("A1:F1").Font.Bold = $true

# Row 2 - Prediction Period Open
$emailTypes.Range("A2").Value = 1
$emailTypes.Range("B2").Value = "EMAILTYPE_PREDICTIONPERIODOPEN"
$emailTypes.Range("C2").Value = "Email sent out just after the prediction period has opened"
$emailTypes.Range("D2").Value = "Prediction Period Open"
$emailTypes.Range("E2").Value = 1
$emailTypes.Range("F2").Formula = '="("&A2&",''"&B2&"'',''"&C2&"'',''"&D2&"'',"&E2&"),"'

# Row 3 - Prediction Period Close
$emailTypes.Range("A3").Value = 2
$emailTypes.Range("B3").Value = "EMAILTYPE_PREDICTIONPERIODCLOSE"
$emailTypes.Range("C3").Value = "Email sent out about a week before the prediction period closes"
$emailTypes.Range("D3").Value = "Prediction Period Close"
$emailTypes.Range("E3").Value = 1
$emailTypes.Range("F3").Formula = '="("&A3&",''"&B3&"'',''"&C3&"'',''"&D3&"'',"&E3&"),"'

# Row 4 - Email Verification (new row)
$emailTypes.Range("A4").Value = 3
$emailTypes.Range("B4").Value = "EMAILTYPE_EMAILVERIFICATION"
$emailTypes.Range("C4").Value = "Email to verify user's email after registration"
$emailTypes.Range("D4").Value = "Email Verification"
$emailTypes.Range("E4").Value = 0
$emailTypes.Range("F4").Formula = '="("&A4&",''"&B4&"'',''"&C4&"'',''"&D4&"'',"&E4&"),"'

# Widen the new "IsSubscription" column roughly to the author's width (~24 chars)
$emailTypes.Columns.Item(5).ColumnWidth = 23.1

$emailTypes.Range("C4").Select() | Out-Null

# ---------------------------------------------------------------------------
# 2. Insert a brand-new "Dates" worksheet right after "EmailTypes" (and
#    before "SQL Commands"), documenting datetime columns / types used
#    across the DB & the Java/SQL type mapping.
# ---------------------------------------------------------------------------
$dates = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $emailTypes)
$dates.Name = "Dates"

# Title
$dates.Range("A1").Value = "Dates"
$dates.Range("A1").Font.Italic = $true

# Header row
$dates.Range("B2").Value = "Table"
$dates.Range("C2").Value = "Column"
$dates.Range("D2").Value = "SQL DT"
$dates.Range("E2").Value = "Java DT"
$dates.Range("B2:E2").Font.Bold = $true

# Data rows
$dates.Range("B3").Value = "ConfirmationTokens"
$dates.Range("C3").Value = "CreatedDateTimeUTC"
$dates.Range("D3").Value = "datetime"
$dates.Range("E3").Value = "Instant"

$dates.Range("B4").Value = "Events"
$dates.Range("C4").Value = "EventDateTimeUTC"
$dates.Range("D4").Value = "datetime"
$dates.Range("E4").Value = "ZonedDateTime"

$dates.Range("B5").Value = "PeriodPredictions"
$dates.Range("C5").Value = "PredictionDateTimeUTC"
$dates.Range("D5").Value = "datetime"
$dates.Range("E5").Value = "Instant"

$dates.Range("B6").Value = "RefreshTokens"
$dates.Range("C6").Value = "LastUsageDateTimeUTC"
$dates.Range("D6").Value = "datetime"
$dates.Range("E6").Value = "Instant"

$dates.Range("B7").Value = "Users"
$dates.Range("C7").Value = "UserCreateDateTimeUTC"
$dates.Range("D7").Value = "datetime"
$dates.Range("E7").Value = "Instant"

$dates.Range("B8").Value = "EmailHistory"
$dates.Range("C8").Value = "RowCreatedDateTimeUTC"
$dates.Range("D8").Value = "datetime"
$dates.Range("E8").Value = "Instant"

$dates.Range("B9").Value = "EmailHistory"
$dates.Range("C9").Value = "EmailSentDateTimeUTC"
$dates.Range("D9").Value = "datetime"
$dates.Range("E9").Value = "Instant"

$dates.Range("B2:C9").EntireColumn.AutoFit() | Out-Null

$dates.PageSetup.PaperSize = 9
$dates.PageSetup.Orientation = 1

$dates.Range("D8").Select() | Out-Null
